# Add the 10 May 2020 ("10 Mayıs 2020") row of data to the covid19-turkey
# "data" worksheet table (Table3).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet has a single table (Table3) covering A1:E59 with a header
# row plus 58 data rows. Add one more row to the table for the new date.
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Fill in the new row's values (row 60: date, test, case, death, recovered).
$ws.Range("A60").Value2 = 43961
$ws.Range("B60").Value2 = 36187
$ws.Range("C60").Value2 = 1542
$ws.Range("D60").Value2 = 47
$ws.Range("E60").Value2 = 3211

# Match the selection recorded in the saved workbook (moves from E58 to E59).
$ws.Range("E59").Select() | Out-Null
